$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be pre-formatted as
# Text so Excel keeps them as strings instead of silently converting them to
# numeric values (matches the original file, where these are text cells).
$textCells = @('D5', 'D6', 'D9', 'D10', 'D12', 'D15', 'D22', 'D24', 'D25', 'D29', 'D31', 'D33', 'D34', 'D35', 'D36', 'D38', 'D40', 'D43', 'D45', 'D48', 'D49', 'D50')
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = '57.506.65'
$ws.Cells.Item(2, 5).Value = '  -0.45%  '
$ws.Cells.Item(3, 4).Value = '3.106.29'
$ws.Cells.Item(3, 5).Value = '  +1.42%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '523.65'
$ws.Cells.Item(5, 5).Value = '  +1.20%  '
$ws.Cells.Item(6, 4).Value = '141.20'
$ws.Cells.Item(6, 5).Value = '  -0.46%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 4).Value = '3.106.33'
$ws.Cells.Item(8, 5).Value = '  +1.45%  '
$ws.Cells.Item(9, 4).Value = '0.436'
$ws.Cells.Item(9, 5).Value = '  +0.12%  '
$ws.Cells.Item(10, 4).Value = '7.24'
$ws.Cells.Item(10, 5).Value = '  -0.82%  '
$ws.Cells.Item(11, 5).Value = '  +0.84%  '
$ws.Cells.Item(12, 4).Value = '0.385'
$ws.Cells.Item(12, 5).Value = '  +1.83%  '
$ws.Cells.Item(13, 4).Value = '3.640.17'
$ws.Cells.Item(13, 5).Value = '  +1.45%  '
$ws.Cells.Item(14, 5).Value = '  +1.07%  '
$ws.Cells.Item(15, 4).Value = '26.14'
$ws.Cells.Item(15, 5).Value = '  +0.74%  '
$ws.Cells.Item(16, 5).Value = '  +0.47%  '
$ws.Cells.Item(17, 4).Value = '57.575.64'
$ws.Cells.Item(17, 5).Value = '  -0.37%  '
$ws.Cells.Item(18, 4).Value = '3.103.27'
$ws.Cells.Item(18, 5).Value = '  +1.21%  '
$ws.Cells.Item(19, 5).Value = '  +0.27%  '
$ws.Cells.Item(21, 5).Value = '  -0.57%  '
$ws.Cells.Item(22, 4).Value = '336.20'
$ws.Cells.Item(22, 5).Value = '  +1.72%  '
$ws.Cells.Item(23, 5).Value = '  +0.10%  '
$ws.Cells.Item(24, 4).Value = '0.512'
$ws.Cells.Item(24, 5).Value = '  +2.62%  '
$ws.Cells.Item(25, 4).Value = '66.53'
$ws.Cells.Item(25, 5).Value = '  +1.25%  '
$ws.Cells.Item(26, 5).Value = '  -0.58%  '
$ws.Cells.Item(27, 5).Value = '  +0.11%  '
$ws.Cells.Item(28, 4).Value = '0.0' + [string][char]0x2083 + '0922'
$ws.Cells.Item(28, 5).Value = '  +2.07%  '
$ws.Cells.Item(29, 4).Value = '6.53'
$ws.Cells.Item(29, 5).Value = '  +2.48%  '
$ws.Cells.Item(30, 5).Value = '  +0.03%  '
$ws.Cells.Item(31, 4).Value = '7.21'
$ws.Cells.Item(31, 5).Value = '  +0.21%  '
$ws.Cells.Item(32, 5).Value = '  +2.20%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '20.92'
$ws.Cells.Item(33, 5).Value = '  +0.94%  '
$ws.Cells.Item(34, 2).Value = 'Fetch.AI'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(34, 4).Value = '1.20'
$ws.Cells.Item(34, 5).Value = '  +0.44%  '
$ws.Cells.Item(35, 4).Value = '157.01'
$ws.Cells.Item(35, 5).Value = '  +1.46%  '
$ws.Cells.Item(36, 4).Value = '4.65'
$ws.Cells.Item(36, 5).Value = '  +3.37%  '
$ws.Cells.Item(37, 5).Value = '  +2.82%  '
$ws.Cells.Item(38, 4).Value = '27.03'
$ws.Cells.Item(38, 5).Value = '  +0.05%  '
$ws.Cells.Item(39, 5).Value = '  +1.31%  '
$ws.Cells.Item(40, 4).Value = '0.0662'
$ws.Cells.Item(40, 5).Value = '  -1.78%  '
$ws.Cells.Item(41, 4).Value = '3.146.51'
$ws.Cells.Item(41, 5).Value = '  +1.34%  '
$ws.Cells.Item(42, 5).Value = '  +0.70%  '
$ws.Cells.Item(43, 4).Value = '0.686'
$ws.Cells.Item(43, 5).Value = '  +4.79%  '
$ws.Cells.Item(44, 5).Value = '  +10.99%  '
$ws.Cells.Item(45, 4).Value = '36.83'
$ws.Cells.Item(45, 5).Value = '  +0.83%  '
$ws.Cells.Item(46, 5).Value = '  -0.05%  '
$ws.Cells.Item(47, 4).Value = '2.298.94'
$ws.Cells.Item(47, 5).Value = '  +1.99%  '
$ws.Cells.Item(48, 4).Value = '0.0260'
$ws.Cells.Item(48, 5).Value = '  +0.58%  '
$ws.Cells.Item(49, 4).Value = '0.978'
$ws.Cells.Item(49, 5).Value = '  +4.86%  '
$ws.Cells.Item(50, 4).Value = '20.78'
$ws.Cells.Item(50, 5).Value = '  +0.91%  '
$ws.Cells.Item(51, 5).Value = '  +2.28%  '
